# Insert a new data row at row 195 of the active sheet (pushes existing
# rows 195-258 down to 196-259) and populate it with the new weekly
# observation for "Ajo" (Femacal de La Calera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 195. This shifts every
# row from 195 downward by one (old 195 -> 196, ..., old 258 -> 259) and
# grows the sheet dimension from A1:R258 to A1:R259.
$ws.Rows.Item(195).Insert()

# Fill in the new row 195 with the new record's data.
$ws.Cells.Item(195, 1).Value = 3
$ws.Cells.Item(195, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(195, 3).Value = "Coquimbo"
$ws.Cells.Item(195, 4).Value = 44524
$ws.Cells.Item(195, 5).Value = 5
$ws.Cells.Item(195, 6).Value = 100112003
$ws.Cells.Item(195, 7).Value = "Ajo"
$ws.Cells.Item(195, 8).Value = "Chino"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 65
$ws.Cells.Item(195, 11).Value = 16000
$ws.Cells.Item(195, 12).Value = 16500
$ws.Cells.Item(195, 13).Value = 16231
$ws.Cells.Item(195, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(195, 15).Value = "China"
$ws.Cells.Item(195, 16).Value = 1623
$ws.Cells.Item(195, 17).Value = 10
$ws.Cells.Item(195, 18).Value = "Hortaliza"
